$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = 0.19
$ws.Range("B16").Select()
